$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14 mirrors the formatting of row 13 (same per-column styles/borders),
# so copy row 13 -> row 14 first to pick up styles, then overwrite the values
# that actually differ for the new lecturer record.
$ws.Range("A13:F13").Copy($ws.Range("A14:F14"))

$ws.Range("A14").Value = 13
$ws.Range("B14").Value = 20011256
$ws.Range("E14").Value = "'0773070597"
$ws.Range("C14").Value = "Nguyễn Thị Nhung"
# D14 (Nữ) and F14 (email) already carried over correctly from the copy.

# New hyperlink on the email cell of the new row, pointing at the same
# mailto target used by the row above it. Adding a hyperlink re-styles the
# cell with Excel's built-in "Hyperlink" look, so reapply row 13's explicit
# formatting afterwards to match how the sheet's other mail cells look.
$ws.Hyperlinks.Add($ws.Range("F14"), "mailto:ngthihanh@gmail.com")
$ws.Range("F13").Copy($ws.Range("F14"))

# Restore selection like the saved workbook shows.
$ws.Range("E17").Select()
